$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Update the title placeholder text
$s.Shapes.Item(1).TextFrame.TextRange.Text = "Lesson #4"

# Update the subtitle placeholder: remove the first paragraph ("Lesson #4")
# and reset the autofit so the line-spacing reduction is cleared.
$subtitle = $s.Shapes.Item(2)
$subtitle.TextFrame.TextRange.Paragraphs(1).Delete()
$subtitle.TextFrame.AutoSize = 2
